# "Creacion de filtros x profesor y modificacion de paginas"
# Adds three new schedule rows (professors MD / SO / AED) below the
# existing Sheet1 table and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new data rows -------------------------------------------------
# Shared-string table must grow in this order: MD, SO, AED, Jueves, Viernes
# -> write column A (new professor codes) before column C (new day names)
# so new strings are interned in that sequence.
$ws.Range("A5").Value = "MD"
$ws.Range("A6").Value = "SO"
$ws.Range("A7").Value = "AED"

$ws.Range("C5").Value = "Jueves"
$ws.Range("C6").Value = "Viernes"
$ws.Range("C7").Value = "Lunes"

$ws.Range("B5").Value = 2
$ws.Range("D5").Value = 0.45833333333333331
$ws.Range("E5").Value = 25

$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 0.83333333333333337
$ws.Range("E6").Value = 10

$ws.Range("B7").Value = 2
$ws.Range("D7").Value = 0.45833333333333331
$ws.Range("E7").Value = 25

# --- formatting ------------------------------------------------------
# Reuse the exact styles already used by the table (centered alignment,
# plus the "h:mm" / "h:mm:ss" time formats on column D) by copying the
# format from the matching existing rows, instead of building new
# number-format/alignment styles from scratch.
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A7:E7").PasteSpecial(-4122)

$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- selection ---------------------------------------------------
$ws.Range("G4").Select()
